$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily row (row 44) with the same style as prior date rows.
$ws.Range("A44").Value = 45993
$ws.Range("A44").NumberFormat = $ws.Range("A43").NumberFormat

$ws.Range("B44").Value = 98
$ws.Range("C44").Value = 111
$ws.Range("D44").Value = 105
